$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J and K columns for rows 1 through 51
$ws.Range("J1:J51").Value = 0.3
$ws.Range("K1:K51").Value = 0.5

# Update the view: scroll position and selection
$win = $excel.ActiveWindow
$win.ScrollRow = 41
$win.ScrollColumn = 1
$ws.Range("K1:K51").Select()
